$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$CENTER = -4108  # xlCenter

# ---------------------------------------------------------------------------
# Seed the five custom cell styles in the exact order the target workbook
# uses them (cellXfs index 1..5), so the emitted styles.xml lines up with
# the authored file:
#   1 = centered
#   2 = centered + h:mm number format
#   3 = centered + hyperlink font
#   4 = centered + bold
#   5 = bold (not centered)
# ---------------------------------------------------------------------------

# 1: centered  (first real use: A1)
$ws.Range("A1").Value = "Пореден Номер"
$ws.Range("A1").HorizontalAlignment = $CENTER

# 2: centered + time format  (first real use: E6)
$ws.Range("E6").Value = "18:45 - 19:00"
$ws.Range("E6").HorizontalAlignment = $CENTER
$ws.Range("E6").NumberFormat = "h:mm"

# 3: centered + hyperlink font  (first real use: D2)
$ws.Range("D2").Value = "p1_21.c"
$ws.Range("D2").HorizontalAlignment = $CENTER
$ws.Hyperlinks.Add($ws.Range("D2"), "p1_21.c")

# 4: centered + bold  (first real use: B1)
$ws.Range("B1").Value = "Задачи"
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").HorizontalAlignment = $CENTER

# 5: bold only  (first real use: D8)
$ws.Range("D8").Value = "Име:"
$ws.Range("D8").Font.Bold = $true

# ---------------------------------------------------------------------------
# Fill in the rest of the table (styles above already exist, so this just
# reuses them).
# ---------------------------------------------------------------------------

# Header row 1 (remaining bold+centered headers)
$ws.Range("C1").Value = "Описание"
$ws.Range("C1").Font.Bold = $true
$ws.Range("C1").HorizontalAlignment = $CENTER

$ws.Range("D1").Value = "Връзка"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").HorizontalAlignment = $CENTER

$ws.Range("E1").Value = "Време"
$ws.Range("E1").Font.Bold = $true
$ws.Range("E1").HorizontalAlignment = $CENTER

# Data rows 2-7
$rows = @(
    @{ n=2; b="Пример1_21";  c="Сумиране на две числа";                                      d="p1_21.c"; e="16:45 - 17:00" },
    @{ n=3; b="Пример2_21";  c="Въвеждане и извеждане на възраст на потребителя";             d="p2_21.c"; e="17:00 - 17:15" },
    @{ n=4; b="Пример3_22";  c="От Целзий в Фаренхайт";                                       d="p3_22.c"; e="17:15 - 17:30" },
    @{ n=5; b="Задача1_22";  c="Отпечатване на трите имена и факултетен номер";               d="z1_22.c"; e="17:45 - 18:30" },
    @{ n=6; b="Задача2_22";  c="самолет - мили в километри ";                                 d="z2_22.c"; e="18:45 - 19:00" },
    @{ n=7; b="Задача3_22";  c="решаване на математическа задача по дадена формула";          d="z3_22.c"; e="19:00 - 19:30" }
)

foreach ($r in $rows) {
    $i = $r.n - 1
    $ws.Range("A$($r.n)").Value = $i
    $ws.Range("A$($r.n)").HorizontalAlignment = $CENTER

    $ws.Range("B$($r.n)").Value = $r.b
    $ws.Range("B$($r.n)").HorizontalAlignment = $CENTER

    $ws.Range("C$($r.n)").Value = $r.c
    $ws.Range("C$($r.n)").HorizontalAlignment = $CENTER

    if ($r.n -ne 2) {
        $ws.Range("D$($r.n)").Value = $r.d
        $ws.Range("D$($r.n)").HorizontalAlignment = $CENTER
        $ws.Hyperlinks.Add($ws.Range("D$($r.n)"), $r.d)
    }

    if ($r.n -ne 6) {
        $ws.Range("E$($r.n)").Value = $r.e
        $ws.Range("E$($r.n)").HorizontalAlignment = $CENTER
    }
}

# Row 8 remainder
$ws.Range("A8").HorizontalAlignment = $CENTER

$ws.Range("E8").Value = "Виктор Мирославов Методиев"
$ws.Range("E8").Font.Bold = $true
$ws.Range("E8").HorizontalAlignment = $CENTER

# ---------------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------------
$ws.Columns("A:B").ColumnWidth = 21.5
$ws.Columns("C").ColumnWidth = 69.66666666666667
$ws.Columns("D").ColumnWidth = 20.833333333333336
$ws.Columns("E").ColumnWidth = 31.833333333333332

# ---------------------------------------------------------------------------
# Page setup
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Selection / view
# ---------------------------------------------------------------------------
$ws.Range("C11").Select() | Out-Null
